$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row rename
$ws.Cells.Item(1,1).Value2 = "mx_state"
$ws.Cells.Item(1,2).Value2 = "mx_municipality"
$ws.Cells.Item(1,3).Value2 = "n_matriculas"
$ws.Cells.Item(1,4).Value2 = "pct_matriculas"

# Title-case corrections for state/municipality names
$ws.Cells.Item(5,2).Value2 = "Rincón De Romos"
$ws.Cells.Item(19,2).Value2 = "Amatenango De La Frontera"
$ws.Cells.Item(22,2).Value2 = "Bejucal De Ocampo"
$ws.Cells.Item(29,2).Value2 = "Comitán De Domínguez"
$ws.Cells.Item(44,2).Value2 = "Mazapa De Madero"
$ws.Cells.Item(46,2).Value2 = "Montecristo De Guerrero"
$ws.Cells.Item(54,2).Value2 = "San Cristóbal De Las Casas"
$ws.Cells.Item(75,2).Value2 = "Guadalupe Y Calvo"
$ws.Cells.Item(94,2).Value2 = "Villa De Álvarez"
$ws.Cells.Item(96,1).Value2 = "Ciudad De México"
$ws.Cells.Item(100,2).Value2 = "Cuajimalpa De Morelos"
$ws.Cells.Item(115,2).Value2 = "Coneto De Comonfort"
$ws.Cells.Item(129,2).Value2 = "San Juan Del Río"
$ws.Cells.Item(135,1).Value2 = "Estado De México"
$ws.Cells.Item(135,2).Value2 = "Acambay De Ruíz Castañeda"
$ws.Cells.Item(137,2).Value2 = "Almoloya De Juárez"
$ws.Cells.Item(141,2).Value2 = "Atizapán De Zaragoza"
$ws.Cells.Item(145,2).Value2 = "Chapa De Mota"
$ws.Cells.Item(148,2).Value2 = "Coacalco De Berriozábal"
$ws.Cells.Item(153,2).Value2 = "Ecatepec De Morelos"
$ws.Cells.Item(157,2).Value2 = "Ixtapan De La Sal"
$ws.Cells.Item(165,2).Value2 = "Naucalpan De Juárez"
$ws.Cells.Item(172,2).Value2 = "San Antonio La Isla"
$ws.Cells.Item(173,2).Value2 = "San Felipe Del Progreso"
$ws.Cells.Item(174,2).Value2 = "San Martín De Las Pirámides"
$ws.Cells.Item(176,2).Value2 = "Soyaniquilpan De Juárez"
$ws.Cells.Item(183,2).Value2 = "Tenango Del Valle"
$ws.Cells.Item(188,2).Value2 = "Tlalnepantla De Baz"
$ws.Cells.Item(193,2).Value2 = "Valle De Bravo"
$ws.Cells.Item(194,2).Value2 = "Villa Del Carbón"
$ws.Cells.Item(206,2).Value2 = "Apaseo El Alto"
$ws.Cells.Item(207,2).Value2 = "Apaseo El Grande"
$ws.Cells.Item(214,2).Value2 = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Cells.Item(217,2).Value2 = "Jaral Del Progreso"
$ws.Cells.Item(228,2).Value2 = "San Diego De La Unión"
$ws.Cells.Item(230,2).Value2 = "San Francisco Del Rincón"
$ws.Cells.Item(232,2).Value2 = "San Luis De La Paz"
$ws.Cells.Item(233,2).Value2 = "Santa Cruz De Juventino Rosas"
$ws.Cells.Item(234,2).Value2 = "Silao De La Victoria"
$ws.Cells.Item(240,2).Value2 = "Acapulco De Juárez"
$ws.Cells.Item(242,2).Value2 = "Ajuchitlán Del Progreso"
$ws.Cells.Item(244,2).Value2 = "Atenango Del Río"
$ws.Cells.Item(246,2).Value2 = "Atoyac De Álvarez"
$ws.Cells.Item(247,2).Value2 = "Ayutla De Los Libres"
$ws.Cells.Item(249,2).Value2 = "Buenavista De Cuéllar"
$ws.Cells.Item(250,2).Value2 = "Chilapa De Álvarez"
$ws.Cells.Item(251,2).Value2 = "Chilpancingo De Los Bravo"
$ws.Cells.Item(253,2).Value2 = "Coyuca De Benítez"
$ws.Cells.Item(254,2).Value2 = "Coyuca De Catalán"
$ws.Cells.Item(257,2).Value2 = "Cuetzala Del Progreso"
$ws.Cells.Item(258,2).Value2 = "Cutzamala De Pinzón"
$ws.Cells.Item(262,2).Value2 = "Huitzuco De Los Figueroa"
$ws.Cells.Item(263,2).Value2 = "Iguala De La Independencia"
$ws.Cells.Item(277,2).Value2 = "Taxco De Alarcón"
$ws.Cells.Item(280,2).Value2 = "Tepecoacuilco De Trujano"
$ws.Cells.Item(282,2).Value2 = "Tixtla De Guerrero"
$ws.Cells.Item(285,2).Value2 = "Tlapa De Comonfort"
$ws.Cells.Item(295,2).Value2 = "Atotonilco El Grande"
$ws.Cells.Item(300,2).Value2 = "Cuautepec De Hinojosa"
$ws.Cells.Item(304,2).Value2 = "Huejutla De Reyes"
$ws.Cells.Item(307,2).Value2 = "Jacala De Ledezma"
$ws.Cells.Item(311,2).Value2 = "Mineral Del Chico"
$ws.Cells.Item(312,2).Value2 = "Molango De Escamilla"
$ws.Cells.Item(314,2).Value2 = "Pachuca De Soto"
$ws.Cells.Item(320,2).Value2 = "Tepehuacán De Guerrero"
$ws.Cells.Item(321,2).Value2 = "Tepeji Del Río De Ocampo"
$ws.Cells.Item(322,2).Value2 = "Tezontepec De Aldama"
$ws.Cells.Item(325,2).Value2 = "Tulancingo De Bravo"
$ws.Cells.Item(327,2).Value2 = "Zacualtipán De Ángeles"
$ws.Cells.Item(332,2).Value2 = "Atotonilco El Alto"
$ws.Cells.Item(333,2).Value2 = "Autlán De Navarro"
$ws.Cells.Item(344,2).Value2 = "Jilotlán De Los Dolores"
$ws.Cells.Item(347,2).Value2 = "Lagos De Moreno"
$ws.Cells.Item(351,2).Value2 = "San Juan De Los Lagos"
$ws.Cells.Item(352,2).Value2 = "San Juanito De Escobedo"
$ws.Cells.Item(353,2).Value2 = "San Miguel El Alto"
$ws.Cells.Item(355,2).Value2 = "Talpa De Allende"
$ws.Cells.Item(358,2).Value2 = "Tlajomulco De Zúñiga"
$ws.Cells.Item(361,2).Value2 = "Unión De Tula"
$ws.Cells.Item(378,2).Value2 = "Coalcomán De Vázquez Pallares"
$ws.Cells.Item(423,2).Value2 = "Tiquicheo De Nicolás Romero"
$ws.Cells.Item(449,2).Value2 = "Tetela Del Volcán"
$ws.Cells.Item(450,2).Value2 = "Tlaltizapán De Zapata"
$ws.Cells.Item(461,2).Value2 = "Santa María Del Oro"
$ws.Cells.Item(471,2).Value2 = "San Nicolás De Los Garza"
$ws.Cells.Item(474,2).Value2 = "Acatlán De Pérez Figueroa"
$ws.Cells.Item(481,2).Value2 = "Coicoyán De Las Flores"
$ws.Cells.Item(482,2).Value2 = "Constancia Del Rosario"
$ws.Cells.Item(483,2).Value2 = "El Barrio De La Soledad"
$ws.Cells.Item(484,2).Value2 = "Fresnillo De Trujano"
$ws.Cells.Item(486,2).Value2 = "Heroica Ciudad De Ejutla De Crespo"
$ws.Cells.Item(487,2).Value2 = "Heroica Ciudad De Tlaxiaco"
$ws.Cells.Item(488,2).Value2 = "Huautla De Jiménez"
$ws.Cells.Item(489,2).Value2 = "Ixtlán De Juárez"
$ws.Cells.Item(490,2).Value2 = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Cells.Item(496,2).Value2 = "Mazatlán Villa De Flores"
$ws.Cells.Item(498,2).Value2 = "Miahuatlán De Porfirio Díaz"
$ws.Cells.Item(499,2).Value2 = "Mixistlán De La Reforma"
$ws.Cells.Item(501,2).Value2 = "Nejapa De Madero"
$ws.Cells.Item(502,2).Value2 = "Oaxaca De Juárez"
$ws.Cells.Item(503,2).Value2 = "Ocotlán De Morelos"
$ws.Cells.Item(504,2).Value2 = "Putla Villa De Guerrero"
$ws.Cells.Item(508,2).Value2 = "San Antonino El Alto"
$ws.Cells.Item(510,2).Value2 = "San Antonio De La Cal"
$ws.Cells.Item(512,2).Value2 = "San Baltazar Yatzachi El Bajo"
$ws.Cells.Item(521,2).Value2 = "San José Del Progreso"
$ws.Cells.Item(523,2).Value2 = "San Juan Bautista Lo De Soto"
$ws.Cells.Item(557,2).Value2 = "San Pedro El Alto"
$ws.Cells.Item(579,2).Value2 = "Santa Inés Del Monte"
$ws.Cells.Item(587,2).Value2 = "Santa María Jalapa Del Marqués"
$ws.Cells.Item(599,2).Value2 = "Santiago Del Río"
$ws.Cells.Item(617,2).Value2 = "Tataltepec De Valdés"
$ws.Cells.Item(618,2).Value2 = "Tezoatlán De Segura Y Luna"
$ws.Cells.Item(619,2).Value2 = "Totontepec Villa De Morelos"
$ws.Cells.Item(621,2).Value2 = "Villa De Tututepec De Melchor Ocampo"
$ws.Cells.Item(623,2).Value2 = "Villa Sola De Vega"
$ws.Cells.Item(624,2).Value2 = "Villa Talea De Castro"
$ws.Cells.Item(625,2).Value2 = "Zimatlán De Álvarez"
$ws.Cells.Item(633,2).Value2 = "Ayotoxco De Guerrero"
$ws.Cells.Item(635,2).Value2 = "Chalchicomula De Sesma"
$ws.Cells.Item(650,2).Value2 = "Izúcar De Matamoros"
$ws.Cells.Item(668,2).Value2 = "Tepanco De López"
$ws.Cells.Item(673,2).Value2 = "Tlacotepec De Benito Juárez"
$ws.Cells.Item(687,2).Value2 = "Amealco De Bonfil"
$ws.Cells.Item(689,2).Value2 = "Cadereyta De Montes"
$ws.Cells.Item(693,2).Value2 = "Jalpan De Serra"
$ws.Cells.Item(694,2).Value2 = "Landa De Matamoros"
$ws.Cells.Item(696,2).Value2 = "Pinal De Amoles"
$ws.Cells.Item(698,2).Value2 = "San Juan Del Río"
$ws.Cells.Item(705,2).Value2 = "Ciudad Del Maíz"
$ws.Cells.Item(708,2).Value2 = "Mexquitic De Carmona"
$ws.Cells.Item(713,2).Value2 = "Soledad De Graciano Sánchez"
$ws.Cells.Item(719,2).Value2 = "Villa De Arista"
$ws.Cells.Item(720,2).Value2 = "Villa De Guadalupe"
$ws.Cells.Item(721,2).Value2 = "Villa De Ramos"
$ws.Cells.Item(770,2).Value2 = "Soto La Marina"
$ws.Cells.Item(778,2).Value2 = "Amaxac De Guerrero"
$ws.Cells.Item(782,2).Value2 = "Contla De Juan Cuamatzi"
$ws.Cells.Item(784,2).Value2 = "Ixtacuixtla De Mariano Matamoros"
$ws.Cells.Item(786,2).Value2 = "Papalotla De Xicohténcatl"
$ws.Cells.Item(788,2).Value2 = "Tetla De La Solidaridad"
$ws.Cells.Item(801,2).Value2 = "Castillo De Teayo"
$ws.Cells.Item(803,2).Value2 = "Cazones De Herrera"
$ws.Cells.Item(811,2).Value2 = "Cosamaloapan De Carpio"
$ws.Cells.Item(820,2).Value2 = "Hueyapan De Ocampo"
$ws.Cells.Item(822,2).Value2 = "Ixhuatlán De Madero"
$ws.Cells.Item(834,2).Value2 = "Martínez De La Torre"
$ws.Cells.Item(843,2).Value2 = "Paso De Ovejas"
$ws.Cells.Item(844,2).Value2 = "Paso Del Macho"
$ws.Cells.Item(847,2).Value2 = "Poza Rica De Hidalgo"
$ws.Cells.Item(852,2).Value2 = "Sayula De Alemán"
$ws.Cells.Item(853,2).Value2 = "Soledad De Doblado"
$ws.Cells.Item(867,2).Value2 = "Tlacotepec De Mejía"
$ws.Cells.Item(879,2).Value2 = "Zontecomatlán De López Y Fuentes"
$ws.Cells.Item(886,2).Value2 = "Nochistlán De Mejía"
$ws.Cells.Item(893,2).Value2 = "Villa De Cos"

# Floating point precision corrections in pct_matriculas column
$ws.Cells.Item(3,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(7,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(11,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(14,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(18,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(44,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(92,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(107,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(117,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(119,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(130,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(134,4).Value2 = 0.00925925925925926
$ws.Cells.Item(139,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(157,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(173,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(186,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(194,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(195,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(227,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(263,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(271,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(279,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(284,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(296,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(306,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(323,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(325,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(344,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(382,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(395,4).Value2 = 0.00925925925925926
$ws.Cells.Item(413,4).Value2 = 0.009485094850948507
$ws.Cells.Item(435,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(438,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(455,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(477,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(533,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(560,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(567,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(577,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(587,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(601,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(608,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(613,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(624,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(661,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(662,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(673,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(674,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(687,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(689,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(705,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(743,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(746,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(756,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(766,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(809,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(818,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(828,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(840,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(851,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(856,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(857,4).Value2 = 0.0009033423667570008
$ws.Cells.Item(867,4).Value2 = 0.0009033423667570008

# Remove trailing metadata rows (901-905) and fix dimension
$ws.Range("A901:D905").ClearContents()
